# Auto-generated edit script: update TPM-derived NATMI metrics for Cntf-Il6ra sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.08830399999999999
$ws.Cells.Item(2, 8).Value = 0.264912
$ws.Cells.Item(2, 9).Value = 0.04372337970871547
$ws.Cells.Item(2, 10).Value = 0.04372337970871546
$ws.Cells.Item(2, 13).Value = 0.9317853333333334
$ws.Cells.Item(2, 14).Value = 2.795356
$ws.Cells.Item(2, 15).Value = 0.1255826100074751
$ws.Cells.Item(2, 16).Value = 0.1255826100074751
$ws.Cells.Item(2, 17).Value = 0.08228037207466667
$ws.Cells.Item(2, 18).Value = 0.740523348672
$ws.Cells.Item(2, 19).Value = 0.005490896142168367
$ws.Cells.Item(2, 20).Value = 0.005490896142168366

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.08830399999999999
$ws.Cells.Item(3, 8).Value = 0.264912
$ws.Cells.Item(3, 9).Value = 0.04372337970871547
$ws.Cells.Item(3, 10).Value = 0.04372337970871546
$ws.Cells.Item(3, 15).Value = 0.7447810673036616
$ws.Cells.Item(3, 16).Value = 0.7447810673036616
$ws.Cells.Item(3, 17).Value = 0.4879725252426667
$ws.Cells.Item(3, 18).Value = 4.391752727184
$ws.Cells.Item(3, 19).Value = 0.03256434540558036
$ws.Cells.Item(3, 20).Value = 0.03256434540558036

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.08830399999999999
$ws.Cells.Item(4, 8).Value = 0.264912
$ws.Cells.Item(4, 9).Value = 0.04372337970871547
$ws.Cells.Item(4, 10).Value = 0.04372337970871546
$ws.Cells.Item(4, 13).Value = 0.9618626666666666
$ws.Cells.Item(4, 15).Value = 0.1296363226888633
$ws.Cells.Item(4, 16).Value = 0.1296363226888633
$ws.Cells.Item(4, 17).Value = 0.08493632091733333
$ws.Cells.Item(4, 18).Value = 0.7644268882559999
$ws.Cells.Item(4, 19).Value = 0.005668138160966736
$ws.Cells.Item(4, 20).Value = 0.005668138160966735

# Row 5
$ws.Cells.Item(5, 9).Value = 0.5310748730197871
$ws.Cells.Item(5, 10).Value = 0.531074873019787
$ws.Cells.Item(5, 13).Value = 0.9317853333333334
$ws.Cells.Item(5, 14).Value = 2.795356
$ws.Cells.Item(5, 15).Value = 0.1255826100074751
$ws.Cells.Item(5, 16).Value = 0.1255826100074751
$ws.Cells.Item(5, 17).Value = 0.9993975406906664
$ws.Cells.Item(5, 18).Value = 8.994577866215998
$ws.Cells.Item(5, 19).Value = 0.0666937686632133
$ws.Cells.Item(5, 20).Value = 0.06669376866321329

# Row 6
$ws.Cells.Item(6, 9).Value = 0.5310748730197871
$ws.Cells.Item(6, 10).Value = 0.531074873019787
$ws.Cells.Item(6, 15).Value = 0.7447810673036616
$ws.Cells.Item(6, 16).Value = 0.7447810673036616
$ws.Cells.Item(6, 19).Value = 0.3955345107458336
$ws.Cells.Item(6, 20).Value = 0.3955345107458335

# Row 7
$ws.Cells.Item(7, 9).Value = 0.5310748730197871
$ws.Cells.Item(7, 10).Value = 0.531074873019787
$ws.Cells.Item(7, 13).Value = 0.9618626666666666
$ws.Cells.Item(7, 15).Value = 0.1296363226888633
$ws.Cells.Item(7, 16).Value = 0.1296363226888633
$ws.Cells.Item(7, 18).Value = 9.284916109367998
$ws.Cells.Item(7, 19).Value = 0.06884659361074022
$ws.Cells.Item(7, 20).Value = 0.06884659361074021

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4252017472714976
$ws.Cells.Item(8, 10).Value = 0.4252017472714976
$ws.Cells.Item(8, 13).Value = 0.9317853333333334
$ws.Cells.Item(8, 14).Value = 2.795356
$ws.Cells.Item(8, 15).Value = 0.1255826100074751
$ws.Cells.Item(8, 16).Value = 0.1255826100074751
$ws.Cells.Item(8, 17).Value = 0.8001613371466667
$ws.Cells.Item(8, 18).Value = 7.201452034320001
$ws.Cells.Item(8, 19).Value = 0.0533979452020935
$ws.Cells.Item(8, 20).Value = 0.05339794520209349

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4252017472714976
$ws.Cells.Item(9, 10).Value = 0.4252017472714976
$ws.Cells.Item(9, 15).Value = 0.7447810673036616
$ws.Cells.Item(9, 16).Value = 0.7447810673036616
$ws.Cells.Item(9, 19).Value = 0.3166822111522478
$ws.Cells.Item(9, 20).Value = 0.3166822111522477

# Row 10
$ws.Cells.Item(10, 9).Value = 0.4252017472714976
$ws.Cells.Item(10, 10).Value = 0.4252017472714976
$ws.Cells.Item(10, 13).Value = 0.9618626666666666
$ws.Cells.Item(10, 15).Value = 0.1296363226888633
$ws.Cells.Item(10, 16).Value = 0.1296363226888633
$ws.Cells.Item(10, 17).Value = 0.8259899463733333
$ws.Cells.Item(10, 18).Value = 7.43390951736
$ws.Cells.Item(10, 19).Value = 0.05512159091715636
$ws.Cells.Item(10, 20).Value = 0.05512159091715636
